$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data points for row 2 (tofino client columns) and their average
$ws.Range("M2").Formula = "=AVERAGE(J2:L2)"
$ws.Range("M2").Font.Bold = $true
$ws.Range("O2").Value = "one client, one switch, tofino"

# New data points for row 4 (tofino measurements) and label
$ws.Range("J4").Value = 0.67838543653488104
$ws.Range("K4").Value = 0.66145831346511796
$ws.Range("O4").Value = "two clients, one switch, tofino"

# Restore the selection left behind by the author's editing session
$ws.Range("J13").Select() | Out-Null
